# Project clean up for the final upload
# Add a new "PET Img/ Img" column (H) with a per-row multiplier and a
# derived column (I) that multiplies the "Img" count (C) by that multiplier.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Column I formulas first (so they don't inherit column H's number format) ---
$ws.Range("I3").Formula = "=C3*H3"
$ws.Range("I4").Formula = "=C4*H4"
$ws.Range("I5").Formula = "=C5*H5"
$ws.Range("I6").Formula = "=C6*H6"
$ws.Range("I7").Formula = "=C7*H7"
$ws.Range("I8").Formula = "=C8*H8"

# --- Column H header ---
$ws.Range("H2").Value = "PET Img/ Img"
$ws.Range("H2").Font.Italic = $true
$ws.Range("H2").HorizontalAlignment = -4108

# --- Column H values (integer number format) ---
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 4
$ws.Range("H5").Value = 4
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 3
$ws.Range("H8").Value = 4
$ws.Range("H3:H8").NumberFormat = "0"

# --- Cosmetic: active selection, matching the saved state ---
$ws.Range("I7").Select()
